# FA50_TestData_PrepareSourceLines_21C.xlsx
#
# The author re-uploaded this test-data workbook and, in doing so, scrubbed
# the hard-coded Oracle Cloud URL / username / password that used to live in
# the "Input_Value" sheet's URL / UserName / Password columns (AA2:AC2),
# along with the hyperlink that pointed at that URL.
#
# Clearing those three cells (and dropping the hyperlink on AA2) is the
# actual content-level edit; once those strings are no longer referenced by
# any cell, Excel's own save path naturally drops them from the shared
# string table, which is what shrinks sharedStrings.xml's count/uniqueCount
# (49/46 -> 46/43) and shifts every other <v> string index down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")
$ws.Activate()

# Drop the hyperlink that lived on AA2 (pointed at the Oracle Cloud URL).
$ws.Range("AA2").Hyperlinks.Delete()

# Scrub the URL / UserName / Password sample values.
$ws.Range("AA2").Value = ""
$ws.Range("AB2").Value = ""
$ws.Range("AC2").Value = ""

# Leave the selection on the cells that were just cleared.
$ws.Range("AA2:AC2").Select()
